$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 13: the full OMS frame string gets its FD-field byte fixed (04 -> 03),
# matching the "LittleEndian" correction applied throughout the sheet.
# ---------------------------------------------------------------------------
$ws.Range("A13").Value = "0411A15A993B441100000000426C244C02FD74000003FD17032C0004913C5E6F0100033933F9310259114202657D0B"

# ---------------------------------------------------------------------------
# Row 19 column A: same 04->03 fix for the standalone field value.
# ---------------------------------------------------------------------------
$ws.Range("A19").Value = "03 FD 17 032C00"

# ---------------------------------------------------------------------------
# Drop the old stray note in F20 ("TU CHYBA JEST BLAD") - replaced below by
# the new decoded-bits table (columns G/N/P/R/T).
# ---------------------------------------------------------------------------
$ws.Range("F20").ClearContents()

# ---------------------------------------------------------------------------
# New "decoded bits" table in column G (one binary string per data row,
# 15-23), stored as text so leading zeros survive.
# ---------------------------------------------------------------------------
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "000001000001000110100001010110101001100100111011"

$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "010001000001000100000000000000000000000000000000"

$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "01000010011011000010010001001100"

$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "0000001011111101011101000000000000000000"

$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "000001001111110100010111000000110010110000000000"

$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "00000100100100010011110001011110011011110000000100000000"

$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "0000001100111001001100111111100100110001"

$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "00000010010110010001000101000010"

$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "00000010011001010111110100001011"

# ---------------------------------------------------------------------------
# Column N: the full concatenated bitstream (row 15 only).
# ---------------------------------------------------------------------------
$ws.Range("N15").NumberFormat = "@"
$ws.Range("N15").Value = "00000100000100011010000101011010100110010011101101000100000100010000000000000000000000000000000001000010011011000010010001001100000000101111110101110100000000000000000000000100111111010001011100000011001011000000000000000100100100010011110001011110011011110000000100000000"

$ws.Range("N20").NumberFormat = "@"
$ws.Range("N20").Value = "00000100"

# ---------------------------------------------------------------------------
# Column P: byte-level breakdown.
# ---------------------------------------------------------------------------
$ws.Range("P18").Value = "03 FD 17 032C00"
$ws.Range("P19").Value = "03 FD 17 032C00"

$ws.Range("P20").NumberFormat = "@"
$ws.Range("P20").Value = "10010001 "

# ---------------------------------------------------------------------------
# Column R: further breakdown.
# ---------------------------------------------------------------------------
$ws.Range("R19").NumberFormat = "@"
$ws.Range("R19").Value = "00010111"

$ws.Range("R20").NumberFormat = "@"
$ws.Range("R20").Value = "00111100"

# ---------------------------------------------------------------------------
# Column T: final breakdown piece (row 20 only).
# ---------------------------------------------------------------------------
$ws.Range("T20").Value = "01011110 01101111 00000001 0000 0000"

# ---------------------------------------------------------------------------
# New rows 26-30: extra worked example block below the table.
# ---------------------------------------------------------------------------
$ws.Range("A26").Value = "04 11 3EB0983B "

$ws.Range("A27").Value = "44 11 00000000"

$ws.Range("A28").Value = "42 6C 244C "

$ws.Range("A29").Value = "02 FD 74 0000 "

$ws.Range("A30").Value = "04FD17032C0004913CC2190200033933F9310259114202651B0C"

# ---------------------------------------------------------------------------
# Column N width (best-fit for the long bitstream column).
# ---------------------------------------------------------------------------
$ws.Columns("N").ColumnWidth = 19.42578125

# ---------------------------------------------------------------------------
# Selection / view state.
# ---------------------------------------------------------------------------
$ws.Range("I14").Select()
